# Update conceptos_df rows: fix retention key assignment logic and remove
# redundant filters. Rows 2-4 get corrected data (service/currency/dates/
# retention code/etc.) and the now-redundant row 5 is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force a literal-text interpretation (mirrors typing `'value` in the
    # Excel UI) so numeric-looking / date-looking strings are not silently
    # coerced into numbers or date serials.
    $range.Value = "'" + $text
}

# ---- Row 2 ----
$ws.Range("F2").Value = "T00003"
$ws.Range("H2").Value = "24-031491"
$ws.Range("I2").Value = "USD"
$ws.Range("K2").Value = "8FCF09CC-58F1-4E6B-ADCA-CA13CD9623CA"
Set-TextValue $ws.Range("M2") "01/11/2024"
Set-TextValue $ws.Range("N2") "01/11/2024"
$ws.Range("P2").Value = "FLETE_TER_N"
Set-TextValue $ws.Range("R2") "75.00"
$ws.Range("T2").Value = "A3V"
Set-TextValue $ws.Range("U2") "12.00"
$ws.Range("V2").Value = "24-031491-290519"
$ws.Range("X2").Value = "Traslado y retención de IVA/ cve SAT: T. local (o con clave general) /Sin CCP"
Set-TextValue $ws.Range("Y2") "78101801"
$ws.Range("Z2").Value = "SERVICIO DE TRANSPORTE DE CARGA LOCAL REF.:24-031491"

# ---- Row 3 ----
$ws.Range("F3").Value = "T00003"
$ws.Range("H3").Value = "24-031489"
$ws.Range("I3").Value = "USD"
$ws.Range("K3").Value = "2E7DB79F-F9F8-42CF-8039-BFC802CF1130"
Set-TextValue $ws.Range("M3") "01/11/2024"
Set-TextValue $ws.Range("N3") "01/11/2024"
$ws.Range("P3").Value = "FLETE_TER_N"
Set-TextValue $ws.Range("R3") "75.00"
$ws.Range("T3").Value = "A3V"
Set-TextValue $ws.Range("U3") "12.00"
$ws.Range("V3").Value = "24-031489-290520"
$ws.Range("X3").Value = "Traslado y retención de IVA/ cve SAT: T. local (o con clave general) /Sin CCP"
Set-TextValue $ws.Range("Y3") "78101801"
$ws.Range("Z3").Value = "SERVICIO DE TRANSPORTE DE CARGA LOCAL REF.:24-031489"

# ---- Row 4 ----
$ws.Range("F4").Value = "T00003"
$ws.Range("H4").Value = "24-031490"
$ws.Range("I4").Value = "USD"
$ws.Range("K4").Value = "CC8C6842-7097-4C01-A2FB-D34397AB06F5"
Set-TextValue $ws.Range("M4") "01/11/2024"
Set-TextValue $ws.Range("N4") "01/11/2024"
$ws.Range("P4").Value = "FLETE_TER_N"
Set-TextValue $ws.Range("R4") "75.00"
$ws.Range("T4").Value = "A3V"
Set-TextValue $ws.Range("U4") "12.00"
$ws.Range("V4").Value = "24-031490-290518"
$ws.Range("X4").Value = "Traslado y retención de IVA/ cve SAT: T. local (o con clave general) /Sin CCP"
Set-TextValue $ws.Range("Y4") "78101801"
$ws.Range("Z4").Value = "SERVICIO DE TRANSPORTE DE CARGA LOCAL REF.:24-031490"

# ---- Row 5 removed entirely (no longer a distinct conceptos row) ----
$ws.Rows.Item(5).Delete()
